$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.263.26"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "1.551.16"
$ws.Range("E3").Value = "  -4.93%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -5.24%  "
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "'17.72"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").Value = "'0.0779"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").Value = "1.766.56"
$ws.Range("D13").Value = "1.548.79"
$ws.Range("E13").Value = "  -6.22%  "
$ws.Range("D14").Value = "'3.97"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("D16").Value = "25.260.28"
$ws.Range("D17").Value = "0.0₃0706"
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("D18").Value = "'58.55"
$ws.Range("E18").Value = "  -4.84%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'185.41"
$ws.Range("E20").Value = "  -4.13%  "
$ws.Range("D21").Value = "'4.09"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'138.62"
$ws.Range("E26").Value = "  -3.81%  "
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("D28").Value = "'14.80"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("E30").Value = "  -6.82%  "
$ws.Range("E31").Value = "  -4.38%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").Value = "1.082.04"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("E39").Value = "  -5.84%  "
$ws.Range("E40").Value = "  -7.84%  "
$ws.Range("E41").Value = "  -10.87%  "
$ws.Range("D42").Value = "'0.796"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("D43").Value = "'92.67"
$ws.Range("E43").Value = "  -5.68%  "
$ws.Range("D44").Value = "'5.04"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "1.682.00"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").Value = "0.0₆0106"
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").Value = "'52.21"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  -2.13%  "
